# Add a new "2022-Q3" sheet (with its fund-holdings detail data) and
# update the "总计" (summary) sheet with a new leading row for it.
#
# Resulting sheet order: 总计, 2022-Q3, 2022-Q2, 2021-Q3, 2021-Q2.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Build the new "2022-Q3" worksheet.
#    Copy the existing "2022-Q2" sheet (same column layout + per-cell
#    styles already correct for 10 rows) so every cell keeps the right
#    formatting, place the copy right before it, rename it, then
#    overwrite the data and drop the now-unneeded trailing rows.
# ------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
# $q2 now resolves (by position) to the freshly inserted copy, which
# landed right before the original sheet - exactly where we want it.
$q3 = $q2
$q3.Name = "2022-Q3"

# Drop the leftover rows (source sheet had 13 data rows, this one only
# needs 9), bringing the used range back down to A1:H10.
$q3.Rows("11:14").Delete()

# Header row (B1:H1) is identical to the source sheet already, so only
# the data rows (2-10) need new values.
#  (the leading "," on each row is required - without it this host's
#   PowerShell interpreter flattens the array-of-arrays into one long
#   flat list instead of keeping each row as its own element)
$q3Data = @(
    ,@(0, "160212", "国泰估值优势混合（LOF）A",             "9.14",  "94.29", "8.15", "0.7449", 3)
    ,@(1, "005535", "泰信竞争优选灵活配置混合",               "11.00", "89.95", "4.75", "0.5225", 7)
    ,@(2, "290006", "泰信蓝筹精选混合",                       "10.65", "90.14", "4.70", "0.5006", 6)
    ,@(3, "020026", "国泰成长优选混合",                       "6.11",  "93.49", "4.65", "0.2841", 8)
    ,@(4, "290002", "泰信先行策略混合",                       "6.17",  "87.99", "4.43", "0.2733", 10)
    ,@(5, "011273", "泰信景气驱动12个月持有期混合A",          "1.09",  "74.12", "4.45", "0.0485", 9)
    ,@(6, "159804", "国寿安保国证创业板中盘精选88ETF",        "1.10",  "98.91", "2.44", "0.0268", 3)
    ,@(7, "011274", "泰信景气驱动12个月持有期混合C",          "0.42",  "74.12", "4.45", "0.0187", 9)
    ,@(8, "016616", "国泰估值优势混合（LOF）C",               "0.00",  "94.29", "8.15", $null,     3)
)

$rowNum = 2
foreach ($row in $q3Data) {
    $q3.Range("A$rowNum").Value = $row[0]
    $q3.Range("B$rowNum").NumberFormat = "@"
    $q3.Range("B$rowNum").Value = $row[1]
    $q3.Range("C$rowNum").NumberFormat = "@"
    $q3.Range("C$rowNum").Value = $row[2]
    $q3.Range("D$rowNum").NumberFormat = "@"
    $q3.Range("D$rowNum").Value = $row[3]
    $q3.Range("E$rowNum").NumberFormat = "@"
    $q3.Range("E$rowNum").Value = $row[4]
    $q3.Range("F$rowNum").NumberFormat = "@"
    $q3.Range("F$rowNum").Value = $row[5]
    if ($null -eq $row[6]) {
        $q3.Range("G$rowNum").NumberFormat = "General"
        $q3.Range("G$rowNum").Value = 0
    } else {
        $q3.Range("G$rowNum").NumberFormat = "@"
        $q3.Range("G$rowNum").Value = $row[6]
    }
    $q3.Range("H$rowNum").Value = $row[7]
    $rowNum = $rowNum + 1
}

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q3 totals as the
#    new first data row, pushing the existing rows down by one.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# New row 5 needs the same style as the existing last row (A column),
# so clone its formatting before filling in values.
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift existing rows down (bottom-up so we don't overwrite data before
# reading it), then write the new 2022-Q3 row on top.
$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q2"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 7.24

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q3"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.05

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 13
$summary.Range("D3").Value = 1.52

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 2.42
